$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row, new Coin (B), new Link (C), new Price (D), new Volume(1h) (E).
# $null means that column is unchanged for that row.
$rows = @(
    @{ Row = 2; B = $null; C = $null; D = "25.826.58"; E = "  +0.13%  " }
    @{ Row = 3; B = $null; C = $null; D = "1.738.34"; E = "  -0.48%  " }
    @{ Row = 4; B = $null; C = $null; D = "1.000"; E = "  -0.08%  " }
    @{ Row = 5; B = $null; C = $null; D = "231.93"; E = "  -1.53%  " }
    @{ Row = 6; B = $null; C = $null; D = $null; E = "  -0.06%  " }
    @{ Row = 7; B = $null; C = $null; D = "0.5172"; E = "  +1.82%  " }
    @{ Row = 8; B = $null; C = $null; D = "0.2803"; E = "  +5.09%  " }
    @{ Row = 9; B = "Dogecoin"; C = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"; D = "0.06113"; E = "  -0.95%  " }
    @{ Row = 10; B = "WrappedEther"; C = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D = "1.749.87"; E = "  +0.24%  " }
    @{ Row = 11; B = "TRON"; C = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; D = "0.07035"; E = "  +1.29%  " }
    @{ Row = 12; B = "Solana"; C = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; D = "15.30"; E = "  +0.26%  " }
    @{ Row = 13; B = "Polygon"; C = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"; D = "0.6456"; E = "  +4.15%  " }
    @{ Row = 14; B = "Polkadot"; C = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"; D = "4.526"; E = "  +1.37%  " }
    @{ Row = 15; B = "Litecoin"; C = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; D = "76.96"; E = "  -1.06%  " }
    @{ Row = 16; B = "BinanceUSD"; C = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"; D = "1.000"; E = "  -0.03%  " }
    @{ Row = 17; B = "Dai"; C = "https://coinranking.com/coin/MoTuySvg7+dai-dai"; D = "1.000"; E = "  -0.08%  " }
    @{ Row = 18; B = "WrappedBTC"; C = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D = "25.823.46"; E = "  +0.00%  " }
    @{ Row = 19; B = "Avalanche"; C = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; D = "11.50"; E = "  -0.63%  " }
    @{ Row = 20; B = "ShibaInu"; C = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; D = "0.000006600"; E = "  -0.48%  " }
    @{ Row = 21; B = "WrappedliquidstakedEther2.0"; C = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; D = "1.974.24"; E = "  +0.46%  " }
    @{ Row = 22; B = "Uniswap"; C = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D = "4.139"; E = "  +2.51%  " }
    @{ Row = 23; B = "Cosmos"; C = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; D = "8.680"; E = "  +5.31%  " }
    @{ Row = 24; B = "Chainlink"; C = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; D = "5.142"; E = "  +0.36%  " }
    @{ Row = 25; B = "Monero"; C = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D = "139.55"; E = "  +2.24%  " }
    @{ Row = 26; B = "Toncoin"; C = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; D = "1.512"; E = "  +3.69%  " }
    @{ Row = 27; B = "EthereumClassic"; C = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D = "15.06"; E = "  +0.39%  " }
    @{ Row = 28; B = "LidoDAOToken"; C = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"; D = "1.812"; E = "  +3.02%  " }
    @{ Row = 29; B = "BitcoinCash"; C = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D = "102.13"; E = "  -0.35%  " }
    @{ Row = 30; B = "Stellar"; C = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D = "0.08334"; E = "  +2.41%  " }
    @{ Row = 31; B = "InternetComputer(DFINITY)"; C = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D = "3.681"; E = "  +0.08%  " }
    @{ Row = 32; B = "Filecoin"; C = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D = "3.426"; E = "  +1.36%  " }
    @{ Row = 33; B = "Hedera"; C = "https://coinranking.com/coin/jad286TjB+hedera-hbar"; D = "0.04494"; E = "  +2.41%  " }
    @{ Row = 34; B = "HuobiToken"; C = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; D = "2.609"; E = "  -1.41%  " }
    @{ Row = 35; B = "ARBITRUM"; C = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; D = "0.9848"; E = "  -0.85%  " }
    @{ Row = 36; B = "ImmutableX"; C = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D = "0.6158"; E = "  +2.68%  " }
    @{ Row = 37; B = "MXToken"; C = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"; D = "2.653"; E = "  +3.15%  " }
    @{ Row = 38; B = "VeChain"; C = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D = "0.01583"; E = "  +1.79%  " }
    @{ Row = 39; B = "RenderToken"; C = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; D = "1.941"; E = "  +0.91%  " }
    @{ Row = 40; B = "PaxDollar"; C = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"; D = "0.9995"; E = "  -0.14%  " }
    @{ Row = 41; B = "Quant"; C = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"; D = "100.71"; E = "  -0.91%  " }
    @{ Row = 42; B = "TheSandbox"; C = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"; D = "0.3844"; E = "  +0.79%  " }
    @{ Row = 43; B = "TrustWalletToken"; C = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"; D = "0.7294"; E = "  -2.14%  " }
    @{ Row = 44; B = "FraxShare"; C = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; D = "4.976"; E = "  +1.76%  " }
    @{ Row = 45; B = "Cronos"; C = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; D = "0.05410"; E = "  -1.45%  " }
    @{ Row = 46; B = "Aptos"; C = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D = "6.279"; E = "  +6.28%  " }
    @{ Row = 47; B = "Algorand"; C = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; D = "0.1121"; E = "  +2.26%  " }
    @{ Row = 48; B = "Aave"; C = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"; D = "53.10"; E = "  +1.28%  " }
    @{ Row = 49; B = "Elrond"; C = "https://coinranking.com/coin/omwkOTglq+elrond-egld"; D = "29.99"; E = "  +0.03%  " }
    @{ Row = 50; B = $null; C = $null; D = "7.656"; E = "  +4.02%  " }
    @{ Row = 51; B = "USDD"; C = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"; D = "1.001"; E = "  -0.46%  " }
)

foreach ($item in $rows) {
    $r = $item.Row
    if ($item.D -ne $null) {
        # Force the Price column to remain plain text so values such as
        # "1.000" or "25.826.58" are not reinterpreted as numbers/dates by Excel
        $ws.Cells.Item($r, 4).NumberFormat = "@"
    }
    if ($item.B -ne $null) { $ws.Range("B$r").Value = $item.B }
    if ($item.C -ne $null) { $ws.Range("C$r").Value = $item.C }
    if ($item.D -ne $null) { $ws.Range("D$r").Value = $item.D }
    if ($item.E -ne $null) { $ws.Range("E$r").Value = $item.E }
}
